$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ME")

# New salt/conductivity readings appended below the existing data (rows 18-23)
$data = @(
    @(44302.395833333336, 0,    367.9, 7.4),
    @(44302.394444444442, 5,    367.7, 7.3),
    @(44302.390972222223, 10,   368,   7.3),
    @(44302.388194444444, 15,   367.8, 7.5),
    @(44302.383333333331, 20,   368.6, 7.3),
    @(44302.379861111112, 23.5, 369.6, 7.3)
)

$r = 18
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Move the active sheet/selection from MO to ME, with a new selected cell
$ws.Activate() | Out-Null
$ws.Range("C36").Select() | Out-Null
